$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1.803092333333333
$ws.Range("H2").Value = 5.409276999999999
$ws.Range("I2").Value = 0.1744886524959502
$ws.Range("J2").Value = 0.1744886524959502
$ws.Range("M2").Value = 0.789222
$ws.Range("N2").Value = 2.367666
$ws.Range("O2").Value = 0.01341929863527565
$ws.Range("P2").Value = 0.01341929863527565
$ws.Range("Q2").Value = 1.423040137498
$ws.Range("R2").Value = 12.807361237482
$ws.Range("S2").Value = 0.002341515336309992
$ws.Range("T2").Value = 0.002341515336309991
$ws.Range("G3").Value = 1.803092333333333
$ws.Range("H3").Value = 5.409276999999999
$ws.Range("I3").Value = 0.1744886524959502
$ws.Range("J3").Value = 0.1744886524959502
$ws.Range("O3").Value = 0.005047365584441773
$ws.Range("P3").Value = 0.005047365584441773
$ws.Range("Q3").Value = 0.5352443529653332
$ws.Range("R3").Value = 4.817199176688
$ws.Range("S3").Value = 0.0008807080194836794
$ws.Range("T3").Value = 0.0008807080194836793
$ws.Range("G4").Value = 1.803092333333333
$ws.Range("H4").Value = 5.409276999999999
$ws.Range("I4").Value = 0.1744886524959502
$ws.Range("J4").Value = 0.1744886524959502
$ws.Range("M4").Value = 57.61405833333333
$ws.Range("N4").Value = 172.842175
$ws.Range("O4").Value = 0.9796232927683105
$ws.Range("P4").Value = 0.9796232927683105
$ws.Range("Q4").Value = 103.8834668730528
$ws.Range("R4").Value = 934.9512018574749
$ws.Range("S4").Value = 0.1709331483087883
$ws.Range("T4").Value = 0.1709331483087882
$ws.Range("G5").Value = 1.803092333333333
$ws.Range("H5").Value = 5.409276999999999
$ws.Range("I5").Value = 0.1744886524959502
$ws.Range("J5").Value = 0.1744886524959502
$ws.Range("M5").Value = 0.1123343333333333
$ws.Range("N5").Value = 0.337003
$ws.Range("O5").Value = 0.001910043011972043
$ws.Range("P5").Value = 0.001910043011972043
$ws.Range("Q5").Value = 0.2025491752034444
$ws.Range("R5").Value = 1.822942576831
$ws.Range("S5").Value = 0.0003332808313683079
$ws.Range("T5").Value = 0.0003332808313683079
$ws.Range("I6").Value = 0.4384883998568034
$ws.Range("J6").Value = 0.4384883998568034
$ws.Range("M6").Value = 0.789222
$ws.Range("N6").Value = 2.367666
$ws.Range("O6").Value = 0.01341929863527565
$ws.Range("P6").Value = 0.01341929863527565
$ws.Range("Q6").Value = 3.576086948336
$ws.Range("R6").Value = 32.184782535024
$ws.Range("S6").Value = 0.005884206785782605
$ws.Range("T6").Value = 0.005884206785782604
$ws.Range("I7").Value = 0.4384883998568034
$ws.Range("J7").Value = 0.4384883998568034
$ws.Range("O7").Value = 0.005047365584441773
$ws.Range("P7").Value = 0.005047365584441773
$ws.Range("S7").Value = 0.002213211258614173
$ws.Range("T7").Value = 0.002213211258614172
$ws.Range("I8").Value = 0.4384883998568034
$ws.Range("J8").Value = 0.4384883998568034
$ws.Range("M8").Value = 57.61405833333333
$ws.Range("N8").Value = 172.842175
$ws.Range("O8").Value = 0.9796232927683105
$ws.Range("P8").Value = 0.9796232927683105
$ws.Range("Q8").Value = 261.0582092826889
$ws.Range("R8").Value = 2349.5238835442
$ws.Range("S8").Value = 0.4295534501084294
$ws.Range("T8").Value = 0.4295534501084293
$ws.Range("I9").Value = 0.4384883998568034
$ws.Range("J9").Value = 0.4384883998568034
$ws.Range("M9").Value = 0.1123343333333333
$ws.Range("N9").Value = 0.337003
$ws.Range("O9").Value = 0.001910043011972043
$ws.Range("P9").Value = 0.001910043011972043
$ws.Range("Q9").Value = 0.5090042387102222
$ws.Range("R9").Value = 4.581038148392
$ws.Range("S9").Value = 0.0008375317039772903
$ws.Range("T9").Value = 0.0008375317039772903
$ws.Range("G10").Value = 3.895605666666667
$ws.Range("H10").Value = 11.686817
$ws.Range("I10").Value = 0.3769851220961256
$ws.Range("J10").Value = 0.3769851220961256
$ws.Range("M10").Value = 0.789222
$ws.Range("N10").Value = 2.367666
$ws.Range("O10").Value = 0.01341929863527565
$ws.Range("P10").Value = 0.01341929863527565
$ws.Range("Q10").Value = 3.074497695458
$ws.Range("R10").Value = 27.670479259122
$ws.Range("S10").Value = 0.005058875934463762
$ws.Range("T10").Value = 0.005058875934463762
$ws.Range("G11").Value = 3.895605666666667
$ws.Range("H11").Value = 11.686817
$ws.Range("I11").Value = 0.3769851220961256
$ws.Range("J11").Value = 0.3769851220961256
$ws.Range("O11").Value = 0.005047365584441773
$ws.Range("P11").Value = 0.005047365584441773
$ws.Range("Q11").Value = 1.156402750938667
$ws.Range("R11").Value = 10.407624758448
$ws.Range("S11").Value = 0.001902781731114564
$ws.Range("T11").Value = 0.001902781731114564
$ws.Range("G12").Value = 3.895605666666667
$ws.Range("H12").Value = 11.686817
$ws.Range("I12").Value = 0.3769851220961256
$ws.Range("J12").Value = 0.3769851220961256
$ws.Range("M12").Value = 57.61405833333333
$ws.Range("N12").Value = 172.842175
$ws.Range("O12").Value = 0.9796232927683105
$ws.Range("P12").Value = 0.9796232927683105
$ws.Range("Q12").Value = 224.4416521229972
$ws.Range("R12").Value = 2019.974869106975
$ws.Range("S12").Value = 0.3693034066324701
$ws.Range("T12").Value = 0.3693034066324701
$ws.Range("G13").Value = 3.895605666666667
$ws.Range("H13").Value = 11.686817
$ws.Range("I13").Value = 0.3769851220961256
$ws.Range("J13").Value = 0.3769851220961256
$ws.Range("M13").Value = 0.1123343333333333
$ws.Range("N13").Value = 0.337003
$ws.Range("O13").Value = 0.001910043011972043
$ws.Range("P13").Value = 0.001910043011972043
$ws.Range("Q13").Value = 0.4376102654945556
$ws.Range("R13").Value = 3.938492389451
$ws.Range("S13").Value = 0.0007200577980771321
$ws.Range("T13").Value = 0.0007200577980771322
$ws.Range("G14").Value = 0.1037266666666667
$ws.Range("H14").Value = 0.31118
$ws.Range("I14").Value = 0.01003782555112075
$ws.Range("J14").Value = 0.01003782555112075
$ws.Range("M14").Value = 0.789222
$ws.Range("N14").Value = 2.367666
$ws.Range("O14").Value = 0.01341929863527565
$ws.Range("P14").Value = 0.01341929863527565
$ws.Range("Q14").Value = 0.08186336732000001
$ws.Range("R14").Value = 0.7367703058799999
$ws.Range("S14").Value = 0.0001347005787192897
$ws.Range("T14").Value = 0.0001347005787192897
$ws.Range("G15").Value = 0.1037266666666667
$ws.Range("H15").Value = 0.31118
$ws.Range("I15").Value = 0.01003782555112075
$ws.Range("J15").Value = 0.01003782555112075
$ws.Range("O15").Value = 0.005047365584441773
$ws.Range("P15").Value = 0.005047365584441773
$ws.Range("Q15").Value = 0.03079105354666667
$ws.Range("R15").Value = 0.27711948192
$ws.Range("S15").Value = 0.00005066457522935716
$ws.Range("T15").Value = 0.00005066457522935715
$ws.Range("G16").Value = 0.1037266666666667
$ws.Range("H16").Value = 0.31118
$ws.Range("I16").Value = 0.01003782555112075
$ws.Range("J16").Value = 0.01003782555112075
$ws.Range("M16").Value = 57.61405833333333
$ws.Range("N16").Value = 172.842175
$ws.Range("O16").Value = 0.9796232927683105
$ws.Range("P16").Value = 0.9796232927683105
$ws.Range("Q16").Value = 5.976114224055556
$ws.Range("R16").Value = 53.7850280165
$ws.Range("S16").Value = 0.009833287718622791
$ws.Range("T16").Value = 0.009833287718622789
$ws.Range("G17").Value = 0.1037266666666667
$ws.Range("H17").Value = 0.31118
$ws.Range("I17").Value = 0.01003782555112075
$ws.Range("J17").Value = 0.01003782555112075
$ws.Range("M17").Value = 0.1123343333333333
$ws.Range("N17").Value = 0.337003
$ws.Range("O17").Value = 0.001910043011972043
$ws.Range("P17").Value = 0.001910043011972043
$ws.Range("Q17").Value = 0.01165206594888889
$ws.Range("R17").Value = 0.10486859354
$ws.Range("S17").Value = 0.00001917267854931261
$ws.Range("T17").Value = 0.00001917267854931261

Write-Output "Updated 174 cells with new TPM values"
